$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A16").Value = 5
$ws.Range("B16").Value = "3；14-7；15"
$ws.Range("C16").Value = "goto指令 关机程序 正式进入函数 一些库函数"

$ws.Range("C16").Select()
